$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("50×17=850", $true, $false, $false, $false, $false, $true, 1, $false, "27×47=1269", 2)
$null = $d.Content.Find.Execute("59×77=4543", $true, $false, $false, $false, $false, $true, 1, $false, "27×25=675", 2)
$null = $d.Content.Find.Execute("97×29=2813", $true, $false, $false, $false, $false, $true, 1, $false, "49×13=637", 2)
$null = $d.Content.Find.Execute("16×93=1488", $true, $false, $false, $false, $false, $true, 1, $false, "84×60=5040", 2)
$null = $d.Content.Find.Execute("50×65=3250", $true, $false, $false, $false, $false, $true, 1, $false, "69×84=5796", 2)
$null = $d.Content.Find.Execute("80×68=5440", $true, $false, $false, $false, $false, $true, 1, $false, "76×64=4864", 2)
$null = $d.Content.Find.Execute("98×23=2254", $true, $false, $false, $false, $false, $true, 1, $false, "58×79=4582", 2)
$null = $d.Content.Find.Execute("98×71=6958", $true, $false, $false, $false, $false, $true, 1, $false, "34×16=544", 2)
$null = $d.Content.Find.Execute("54×97=5238", $true, $false, $false, $false, $false, $true, 1, $false, "79×23=1817", 2)
$null = $d.Content.Find.Execute("80×18=1440", $true, $false, $false, $false, $false, $true, 1, $false, "89×31=2759", 2)
$null = $d.Content.Find.Execute("58×24=1392", $true, $false, $false, $false, $false, $true, 1, $false, "62×71=4402", 2)
$null = $d.Content.Find.Execute("56×15=840", $true, $false, $false, $false, $false, $true, 1, $false, "24×22=528", 2)
$null = $d.Content.Find.Execute("23×100=2300", $true, $false, $false, $false, $false, $true, 1, $false, "96×62=5952", 2)
$null = $d.Content.Find.Execute("70×12=840", $true, $false, $false, $false, $false, $true, 1, $false, "61×21=1281", 2)
$null = $d.Content.Find.Execute("54×58=3132", $true, $false, $false, $false, $false, $true, 1, $false, "53×28=1484", 2)
$null = $d.Content.Find.Execute("87×61=5307", $true, $false, $false, $false, $false, $true, 1, $false, "46×45=2070", 2)
$null = $d.Content.Find.Execute("55×83=4565", $true, $false, $false, $false, $false, $true, 1, $false, "79×42=3318", 2)
$null = $d.Content.Find.Execute("99×12=1188", $true, $false, $false, $false, $false, $true, 1, $false, "61×98=5978", 2)
$null = $d.Content.Find.Execute("13×57=741", $true, $false, $false, $false, $false, $true, 1, $false, "91×61=5551", 2)
$null = $d.Content.Find.Execute("82×52=4264", $true, $false, $false, $false, $false, $true, 1, $false, "84×91=7644", 2)
$null = $d.Content.Find.Execute("51×78=3978", $true, $false, $false, $false, $false, $true, 1, $false, "36×63=2268", 2)
$null = $d.Content.Find.Execute("66×86=5676", $true, $false, $false, $false, $false, $true, 1, $false, "83×56=4648", 2)
$null = $d.Content.Find.Execute("21×59=1239", $true, $false, $false, $false, $false, $true, 1, $false, "72×29=2088", 2)
$null = $d.Content.Find.Execute("92×31=2852", $true, $false, $false, $false, $false, $true, 1, $false, "58×21=1218", 2)
$null = $d.Content.Find.Execute("89×74=6586", $true, $false, $false, $false, $false, $true, 1, $false, "31×23=713", 2)
$null = $d.Content.Find.Execute("26×56=1456", $true, $false, $false, $false, $false, $true, 1, $false, "84×84=7056", 2)
$null = $d.Content.Find.Execute("50×28=1400", $true, $false, $false, $false, $false, $true, 1, $false, "29×73=2117", 2)
$null = $d.Content.Find.Execute("27×11=297", $true, $false, $false, $false, $false, $true, 1, $false, "28×48=1344", 2)
$null = $d.Content.Find.Execute("21×31=651", $true, $false, $false, $false, $false, $true, 1, $false, "79×71=5609", 2)
$null = $d.Content.Find.Execute("88×86=7568", $true, $false, $false, $false, $false, $true, 1, $false, "16×18=288", 2)
$null = $d.Content.Find.Execute("30×19=570", $true, $false, $false, $false, $false, $true, 1, $false, "12×15=180", 2)
$null = $d.Content.Find.Execute("89×27=2403", $true, $false, $false, $false, $false, $true, 1, $false, "13×30=390", 2)
$null = $d.Content.Find.Execute("98×11=1078", $true, $false, $false, $false, $false, $true, 1, $false, "90×56=5040", 2)
$null = $d.Content.Find.Execute("32×23=736", $true, $false, $false, $false, $false, $true, 1, $false, "12×80=960", 2)
$null = $d.Content.Find.Execute("14×90=1260", $true, $false, $false, $false, $false, $true, 1, $false, "29×20=580", 2)
$null = $d.Content.Find.Execute("51×41=2091", $true, $false, $false, $false, $false, $true, 1, $false, "87×96=8352", 2)
$null = $d.Content.Find.Execute("57×44=2508", $true, $false, $false, $false, $false, $true, 1, $false, "84×15=1260", 2)
$null = $d.Content.Find.Execute("34×68=2312", $true, $false, $false, $false, $false, $true, 1, $false, "98×46=4508", 2)
$null = $d.Content.Find.Execute("34×26=884", $true, $false, $false, $false, $false, $true, 1, $false, "37×79=2923", 2)
$null = $d.Content.Find.Execute("53×86=4558", $true, $false, $false, $false, $false, $true, 1, $false, "64×90=5760", 2)
$null = $d.Content.Find.Execute("84×33=2772", $true, $false, $false, $false, $false, $true, 1, $false, "52×11=572", 2)
$null = $d.Content.Find.Execute("64×45=2880", $true, $false, $false, $false, $false, $true, 1, $false, "50×61=3050", 2)
$null = $d.Content.Find.Execute("88×92=8096", $true, $false, $false, $false, $false, $true, 1, $false, "80×72=5760", 2)
$null = $d.Content.Find.Execute("56×14=784", $true, $false, $false, $false, $false, $true, 1, $false, "57×27=1539", 2)
$null = $d.Content.Find.Execute("79×18=1422", $true, $false, $false, $false, $false, $true, 1, $false, "26×72=1872", 2)
$null = $d.Content.Find.Execute("18×98=1764", $true, $false, $false, $false, $false, $true, 1, $false, "94×62=5828", 2)
$null = $d.Content.Find.Execute("72×78=5616", $true, $false, $false, $false, $false, $true, 1, $false, "24×20=480", 2)
$null = $d.Content.Find.Execute("39×74=2886", $true, $false, $false, $false, $false, $true, 1, $false, "59×32=1888", 2)
$null = $d.Content.Find.Execute("43×59=2537", $true, $false, $false, $false, $false, $true, 1, $false, "38×46=1748", 2)
$null = $d.Content.Find.Execute("48×41=1968", $true, $false, $false, $false, $false, $true, 1, $false, "25×52=1300", 2)
$null = $d.Content.Find.Execute("82×22=1804", $true, $false, $false, $false, $false, $true, 1, $false, "41×78=3198", 2)
$null = $d.Content.Find.Execute("35×20=700", $true, $false, $false, $false, $false, $true, 1, $false, "67×35=2345", 2)
$null = $d.Content.Find.Execute("36×32=1152", $true, $false, $false, $false, $false, $true, 1, $false, "18×78=1404", 2)
$null = $d.Content.Find.Execute("11×43=473", $true, $false, $false, $false, $false, $true, 1, $false, "40×90=3600", 2)
$null = $d.Content.Find.Execute("29×81=2349", $true, $false, $false, $false, $false, $true, 1, $false, "50×85=4250", 2)
$null = $d.Content.Find.Execute("40×41=1640", $true, $false, $false, $false, $false, $true, 1, $false, "43×44=1892", 2)
$null = $d.Content.Find.Execute("69×90=6210", $true, $false, $false, $false, $false, $true, 1, $false, "65×92=5980", 2)
$null = $d.Content.Find.Execute("88×58=5104", $true, $false, $false, $false, $false, $true, 1, $false, "41×15=615", 2)
$null = $d.Content.Find.Execute("39×45=1755", $true, $false, $false, $false, $false, $true, 1, $false, "55×12=660", 2)
$null = $d.Content.Find.Execute("15×29=435", $true, $false, $false, $false, $false, $true, 1, $false, "90×99=8910", 2)
$null = $d.Content.Find.Execute("92×97=8924", $true, $false, $false, $false, $false, $true, 1, $false, "91×85=7735", 2)
$null = $d.Content.Find.Execute("35×19=665", $true, $false, $false, $false, $false, $true, 1, $false, "32×94=3008", 2)
$null = $d.Content.Find.Execute("85×57=4845", $true, $false, $false, $false, $false, $true, 1, $false, "50×31=1550", 2)
$null = $d.Content.Find.Execute("94×63=5922", $true, $false, $false, $false, $false, $true, 1, $false, "51×39=1989", 2)
$null = $d.Content.Find.Execute("75×61=4575", $true, $false, $false, $false, $false, $true, 1, $false, "22×44=968", 2)
$null = $d.Content.Find.Execute("82×89=7298", $true, $false, $false, $false, $false, $true, 1, $false, "50×90=4500", 2)
$null = $d.Content.Find.Execute("51×52=2652", $true, $false, $false, $false, $false, $true, 1, $false, "31×12=372", 2)
$null = $d.Content.Find.Execute("98×84=8232", $true, $false, $false, $false, $false, $true, 1, $false, "87×58=5046", 2)
$null = $d.Content.Find.Execute("42×15=630", $true, $false, $false, $false, $false, $true, 1, $false, "28×43=1204", 2)
$null = $d.Content.Find.Execute("42×24=1008", $true, $false, $false, $false, $false, $true, 1, $false, "74×29=2146", 2)
$null = $d.Content.Find.Execute("40×14=560", $true, $false, $false, $false, $false, $true, 1, $false, "58×78=4524", 2)
$null = $d.Content.Find.Execute("28×51=1428", $true, $false, $false, $false, $false, $true, 1, $false, "33×57=1881", 2)
$null = $d.Content.Find.Execute("89×77=6853", $true, $false, $false, $false, $false, $true, 1, $false, "60×11=660", 2)
$null = $d.Content.Find.Execute("40×83=3320", $true, $false, $false, $false, $false, $true, 1, $false, "48×21=1008", 2)
$null = $d.Content.Find.Execute("46×62=2852", $true, $false, $false, $false, $false, $true, 1, $false, "94×70=6580", 2)
$null = $d.Content.Find.Execute("90×22=1980", $true, $false, $false, $false, $false, $true, 1, $false, "58×65=3770", 2)
$null = $d.Content.Find.Execute("30×62=1860", $true, $false, $false, $false, $false, $true, 1, $false, "64×57=3648", 2)
$null = $d.Content.Find.Execute("98×54=5292", $true, $false, $false, $false, $false, $true, 1, $false, "19×13=247", 2)
$null = $d.Content.Find.Execute("22×99=2178", $true, $false, $false, $false, $false, $true, 1, $false, "56×87=4872", 2)
$null = $d.Content.Find.Execute("58×70=4060", $true, $false, $false, $false, $false, $true, 1, $false, "21×73=1533", 2)
$null = $d.Content.Find.Execute("81×37=2997", $true, $false, $false, $false, $false, $true, 1, $false, "67×47=3149", 2)
$null = $d.Content.Find.Execute("86×46=3956", $true, $false, $false, $false, $false, $true, 1, $false, "45×85=3825", 2)
$null = $d.Content.Find.Execute("45×71=3195", $true, $false, $false, $false, $false, $true, 1, $false, "53×19=1007", 2)
$null = $d.Content.Find.Execute("31×65=2015", $true, $false, $false, $false, $false, $true, 1, $false, "64×90=5760", 2)
$null = $d.Content.Find.Execute("64×41=2624", $true, $false, $false, $false, $false, $true, 1, $false, "90×49=4410", 2)
$null = $d.Content.Find.Execute("14×15=210", $true, $false, $false, $false, $false, $true, 1, $false, "63×86=5418", 2)
$null = $d.Content.Find.Execute("16×88=1408", $true, $false, $false, $false, $false, $true, 1, $false, "79×69=5451", 2)
$null = $d.Content.Find.Execute("43×65=2795", $true, $false, $false, $false, $false, $true, 1, $false, "61×44=2684", 2)
$null = $d.Content.Find.Execute("83×46=3818", $true, $false, $false, $false, $false, $true, 1, $false, "10×33=330", 2)
$null = $d.Content.Find.Execute("78×89=6942", $true, $false, $false, $false, $false, $true, 1, $false, "93×95=8835", 2)
$null = $d.Content.Find.Execute("95×64=6080", $true, $false, $false, $false, $false, $true, 1, $false, "76×15=1140", 2)
$null = $d.Content.Find.Execute("29×39=1131", $true, $false, $false, $false, $false, $true, 1, $false, "22×21=462", 2)
$null = $d.Content.Find.Execute("35×49=1715", $true, $false, $false, $false, $false, $true, 1, $false, "66×69=4554", 2)
$null = $d.Content.Find.Execute("16×56=896", $true, $false, $false, $false, $false, $true, 1, $false, "65×25=1625", 2)
$null = $d.Content.Find.Execute("99×68=6732", $true, $false, $false, $false, $false, $true, 1, $false, "28×16=448", 2)
$null = $d.Content.Find.Execute("80×61=4880", $true, $false, $false, $false, $false, $true, 1, $false, "50×69=3450", 2)
$null = $d.Content.Find.Execute("84×41=3444", $true, $false, $false, $false, $false, $true, 1, $false, "41×54=2214", 2)
$null = $d.Content.Find.Execute("71×81=5751", $true, $false, $false, $false, $false, $true, 1, $false, "28×26=728", 2)
$null = $d.Content.Find.Execute("29×62=1798", $true, $false, $false, $false, $false, $true, 1, $false, "15×76=1140", 2)
$null = $d.Content.Find.Execute("97×81=7857", $true, $false, $false, $false, $false, $true, 1, $false, "92×10=920", 2)
